# Weekly price update ("Fruta / hortaliza, semanal"): a new price record for
# Espinaca (Mercado Mayorista Lo Valledor de Santiago) is inserted right
# before the existing row 484, pushing the existing rows 484:601 down to
# 485:602 (dimension grows from A1:R601 to A1:R602).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 484, shifting rows 484:601 down to 485:602.
$ws.Rows("484:484").Insert()

# Populate the newly inserted row 484 with the new data record.
$ws.Range("A484").Value = 6
$ws.Range("B484").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C484").Value = "Metropolitana"
$ws.Range("D484").Value = 44798
$ws.Range("E484").Value = 13
$ws.Range("F484").Value = 100112012
$ws.Range("G484").Value = "Espinaca"
$ws.Range("H484").Value = "Sin especificar"
$ws.Range("I484").Value = "Primera"
$ws.Range("J484").Value = 560
$ws.Range("K484").Value = 5500
$ws.Range("L484").Value = 6000
$ws.Range("M484").Value = 5688
$ws.Range("N484").Value = "`$/cuna 10 kilos"
$ws.Range("O484").Value = "Región Metropolitana"
$ws.Range("P484").Value = 569
$ws.Range("Q484").Value = 10
$ws.Range("R484").Value = "Hortaliza"
